$wb = $excel.ActiveWorkbook

# --- Insert the new "ImportFiledata" worksheet after "LogOutData" ---
$afterSheet = $wb.Worksheets.Item("LogOutData")
$importSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$importSheet.Name = "ImportFiledata"

# Header row
$importSheet.Range("A1").Value = "username"
$importSheet.Range("B1").Value = "password"
$importSheet.Range("C1").Value = "test_file_folder"
$importSheet.Range("D1").Value = "test_file_name"

# Data row (first two columns reference the Common sheet, same as other data sheets)
$importSheet.Range("A2").Formula = "=Common!B7"
$importSheet.Range("B2").Formula = "=Common!B5"
$importSheet.Range("C2").Value = "TestData"
$importSheet.Range("D2").Value = "meals.csv"

# Column widths (best-fit sizing similar to the other generated sheets)
$importSheet.Columns.Item(1).ColumnWidth = 17.85
$importSheet.Columns.Item(2).ColumnWidth = 24.675
$importSheet.Columns.Item(3).ColumnWidth = 17.025
$importSheet.Columns.Item(4).ColumnWidth = 24.675

# Selection/active cell on the new sheet
$importSheet.Range("B3").Select() | Out-Null

# --- Update selection on the Common sheet ---
$common = $wb.Worksheets.Item("Common")
$common.Range("B5").Select() | Out-Null

# --- Update selection on the SigninData sheet (it is no longer the active tab) ---
$signinData = $wb.Worksheets.Item("SigninData")
$signinData.Range("A2").Select() | Out-Null

# Make sure the newly inserted sheet ends up as the active tab
$importSheet.Activate() | Out-Null
$importSheet.Range("B3").Select() | Out-Null
